{"js": "// The captured diff for this document is a pure XML re-serialization\n// artifact: every hunk only reorders the attributes that already exist\n// on an element (namespace declarations on <w:document>, and attributes\n// such as w:color/@w:val/@w:themeColor/@w:themeShade, w:pgSz/@w:w/@w:h,\n// w:pgMar/*, w:rFonts/*, w:lang/*, w:latentStyles/*, w:lsdException/*,\n// w:style/*, w:tblInd/*, w:tblCellMar/* ...). No text, formatting value,\n// style, or structural content actually changes -- the commit message\n// (\"Fixed POI packaging and upgraded to POI 3.15\") confirms this was a\n// tooling/serializer change, not a content edit.\n//\n// The Word JavaScript API only exposes the document's semantic object\n// model (text, ranges, styles, properties, ...); it has no notion of\n// \"attribute order\" inside a start tag, which is not semantically\n// meaningful in XML/OOXML in the first place. So the correct, faithful\n// replay of this change through Office.js is to leave the document's\n// content and formatting exactly as they are.\n//\n// Touch the body in a read-only way so the script still exercises the\n// context/sync pattern requested, without mutating anything.\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n", "ps1": "# The captured diff for this document is a pure XML re-serialization\n# artifact: every hunk only reorders attributes that already exist on an\n# element (namespace declarations on <w:document>, and attributes such\n# as w:color/@w:val/@w:themeColor/@w:themeShade, w:pgSz/@w:w/@w:h,\n# w:pgMar/*, w:rFonts/*, w:lang/*, w:latentStyles/*, w:lsdException/*,\n# w:style/*, w:tblInd/*, w:tblCellMar/* ...). No text, formatting value,\n# style, or structural content actually changes -- the commit message\n# (\"Fixed POI packaging and upgraded to POI 3.15\") confirms this was a\n# tooling/serializer change, not a content edit.\n#\n# The Word COM object model only exposes the document's semantic object\n# model (text, ranges, styles, properties, ...); it has no notion of\n# \"attribute order\" inside a start tag, which is not semantically\n# meaningful in XML/OOXML in the first place. So the correct, faithful\n# replay of this change through COM is to leave the document's content\n# and formatting exactly as they are.\n#\n# Touch the document in a read-only way so the script still exercises\n# the object model as requested, without mutating anything.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
